$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.091.80"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "3.900.55"
$ws.Range("E3").Value = "  +3.32%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'465.20"
$ws.Range("E5").Value = "  +8.72%  "
$ws.Range("D6").Value = "'144.80"
$ws.Range("E6").Value = "  +3.79%  "
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  +7.83%  "
$ws.Range("D11").Value = "'0.0000341"
$ws.Range("E11").Value = "  +8.95%  "
$ws.Range("D12").Value = "'42.94"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").Value = "4.525.41"
$ws.Range("D15").Value = "'15.33"
$ws.Range("E15").Value = "  +2.67%  "
$ws.Range("D16").Value = "3.905.78"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "'19.99"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  +2.06%  "
$ws.Range("D20").Value = "67.314.89"
$ws.Range("E20").Value = "  +1.85%  "
$ws.Range("D21").Value = "'431.67"
$ws.Range("E21").Value = "  +5.82%  "
$ws.Range("D22").Value = "'14.71"
$ws.Range("E22").Value = "  -3.76%  "
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("D24").Value = "'88.63"
$ws.Range("E24").Value = "  +4.11%  "
$ws.Range("D25").Value = "'38.49"
$ws.Range("E25").Value = "  +4.64%  "
$ws.Range("D26").Value = "'3.52"
$ws.Range("E26").Value = "  +7.04%  "
$ws.Range("D27").Value = "'5.72"
$ws.Range("E27").Value = "  +5.58%  "
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("D29").Value = "'9.63"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("D30").Value = "'738.63"
$ws.Range("E30").Value = "  +4.96%  "
$ws.Range("D31").Value = "'13.64"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("D32").Value = "'0.131"
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").Value = "'2.79"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "'42.98"
$ws.Range("E34").Value = "  +4.82%  "
$ws.Range("E35").Value = "  +5.27%  "
$ws.Range("D36").Value = "'58.13"
$ws.Range("E36").Value = "  +2.87%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "0.0₃0785"
$ws.Range("E38").Value = "  +16.83%  "
$ws.Range("D39").Value = "'5.40"
$ws.Range("E39").Value = "  -6.06%  "
$ws.Range("E40").Value = "  +13.75%  "
$ws.Range("D41").Value = "'0.0475"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.140"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "'0.333"
$ws.Range("E44").Value = "  +4.92%  "
$ws.Range("E45").Value = "  +4.81%  "
$ws.Range("D46").Value = "'2.16"
$ws.Range("E46").Value = "  +5.03%  "
$ws.Range("D47").Value = "'3.40"
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("E48").Value = "  -4.36%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'2.90"
$ws.Range("E49").Value = "  +3.09%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "'3.12"
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").Value = "'143.37"
$ws.Range("E51").Value = "  +0.72%  "
